$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '72.334.96'
Set-TextValue $ws.Range("E2") '  +0.00%  '
Set-TextValue $ws.Range("D3") '2.643.44'
Set-TextValue $ws.Range("E3") '  -0.80%  '
Set-TextValue $ws.Range("E4") '  -0.08%  '
Set-TextValue $ws.Range("D5") '589.30'
Set-TextValue $ws.Range("E5") '  -2.15%  '
Set-TextValue $ws.Range("D6") '174.26'
Set-TextValue $ws.Range("E6") '  -2.36%  '
Set-TextValue $ws.Range("E7") '  -0.08%  '
Set-TextValue $ws.Range("D8") '0.519'
Set-TextValue $ws.Range("E8") '  -0.97%  '
Set-TextValue $ws.Range("D9") '0.173'
Set-TextValue $ws.Range("E9") '  -0.47%  '
Set-TextValue $ws.Range("D10") '2.642.01'
Set-TextValue $ws.Range("E10") '  -0.82%  '
Set-TextValue $ws.Range("D11") '0.170'
Set-TextValue $ws.Range("E11") '  +1.08%  '
Set-TextValue $ws.Range("D12") '0.355'
Set-TextValue $ws.Range("E12") '  -0.60%  '
Set-TextValue $ws.Range("D13") '4.94'
Set-TextValue $ws.Range("E13") '  -1.68%  '
Set-TextValue $ws.Range("D14") '3.120.13'
Set-TextValue $ws.Range("E14") '  -1.07%  '
Set-TextValue $ws.Range("D15") '0.0000186'
Set-TextValue $ws.Range("E15") '  -1.20%  '
Set-TextValue $ws.Range("D16") '72.180.07'
Set-TextValue $ws.Range("D17") '25.90'
Set-TextValue $ws.Range("E17") '  -2.67%  '
Set-TextValue $ws.Range("D18") '2.663.77'
Set-TextValue $ws.Range("E18") '  -0.10%  '
Set-TextValue $ws.Range("D19") '12.13'
Set-TextValue $ws.Range("E19") '  +1.14%  '
Set-TextValue $ws.Range("D20") '7.98'
Set-TextValue $ws.Range("E20") '  -0.78%  '
Set-TextValue $ws.Range("D21") '371.88'
Set-TextValue $ws.Range("E21") '  -1.91%  '
Set-TextValue $ws.Range("D22") '4.15'
Set-TextValue $ws.Range("E22") '  -1.33%  '
Set-TextValue $ws.Range("D23") '2.05'
Set-TextValue $ws.Range("E23") '  -0.58%  '
Set-TextValue $ws.Range("E24") '  -0.05%  '
Set-TextValue $ws.Range("D25") '70.97'
Set-TextValue $ws.Range("E25") '  -2.23%  '
Set-TextValue $ws.Range("D26") '4.25'
Set-TextValue $ws.Range("E26") '  -3.48%  '
Set-TextValue $ws.Range("D27") '9.64'
Set-TextValue $ws.Range("E27") '  -3.47%  '
Set-TextValue $ws.Range("D28") '2.776.32'
Set-TextValue $ws.Range("E28") '  -1.10%  '
Set-TextValue $ws.Range("D29") '0.997'
Set-TextValue $ws.Range("E29") '  -0.37%  '
Set-TextValue $ws.Range("D30") '0.0₃0954'
Set-TextValue $ws.Range("E30") '  +0.41%  '
Set-TextValue $ws.Range("D31") '7.99'
Set-TextValue $ws.Range("E31") '  -2.92%  '
Set-TextValue $ws.Range("D32") '497.42'
Set-TextValue $ws.Range("E32") '  -4.91%  '
Set-TextValue $ws.Range("E33") '  -2.33%  '
Set-TextValue $ws.Range("E34") '  -1.18%  '
Set-TextValue $ws.Range("E35") '  -0.05%  '
Set-TextValue $ws.Range("D36") '161.45'
Set-TextValue $ws.Range("E36") '  -1.18%  '
Set-TextValue $ws.Range("D37") '19.27'
Set-TextValue $ws.Range("E37") '  -1.46%  '
Set-TextValue $ws.Range("D38") '0.113'
Set-TextValue $ws.Range("E38") '  +1.72%  '
Set-TextValue $ws.Range("E39") '  -1.27%  '
Set-TextValue $ws.Range("D40") '1.36'
Set-TextValue $ws.Range("E40") '  -2.81%  '
Set-TextValue $ws.Range("E41") '  -0.11%  '
Set-TextValue $ws.Range("E42") '  -5.87%  '
Set-TextValue $ws.Range("D43") '2.57'
Set-TextValue $ws.Range("E43") '  -1.19%  '
Set-TextValue $ws.Range("D44") '4.90'
Set-TextValue $ws.Range("E44") '  -3.31%  '
Set-TextValue $ws.Range("D45") '0.327'
Set-TextValue $ws.Range("E45") '  -2.49%  '
Set-TextValue $ws.Range("D46") '39.05'
Set-TextValue $ws.Range("E46") '  -0.72%  '
Set-TextValue $ws.Range("D47") '152.65'
Set-TextValue $ws.Range("E47") '  -0.30%  '
Set-TextValue $ws.Range("D48") '3.66'
Set-TextValue $ws.Range("E48") '  -2.36%  '
Set-TextValue $ws.Range("D49") '0.548'
Set-TextValue $ws.Range("E49") '  -0.62%  '
Set-TextValue $ws.Range("E50") '  -2.71%  '
Set-TextValue $ws.Range("E51") '  -2.32%  '
